$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of component data appended below the header row
$ws.Range("A2").Value = "OCEBB0015S"
$ws.Range("B2").Value = "Bobine, Modo Comum, 1mH, 30%,  0.8A, 9.2x6x5mm"
$ws.Range("C2").Value = "Prazo de entrega superior a 3 dias"

# Resize the columns to fit the new content
$ws.Range("A1:C2").EntireColumn.AutoFit()

# Leave the selection where the author last left it
$ws.Range("C9").Select()
